$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Здравоохранение" (Healthcare) category column (E) to the
# fourth mini-table (rows 17-20), mirroring the style of the other
# category tables on the sheet (column D holds the matching formatting
# for each row of this table).

$ws.Range("D17").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = "Здравоохранение"

$ws.Range("D18").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "Лечебно-проф. орг. - hospitals (шт.) (8018000)"

$ws.Range("D19").Copy()
$ws.Range("E19").PasteSpecial(-4122)

$ws.Range("D20").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection from E25 to E26, matching the author's cursor
# position after the edit.
$ws.Range("E26").Select()
